
# Update fleet_key, add code 820
# Adds a new row (63) to Sheet1: metier "MIS_C", gear "FPO Burar, tinor, mjärdar", gear_code 820

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 63

$ws.Cells.Item($newRow, 1).Value = "MIS_C"
$ws.Cells.Item($newRow, 2).Value = "FPO Burar, tinor, mjärdar"
$ws.Cells.Item($newRow, 3).Value = 820

# Match the selection left behind by the author (cursor parked on the new row)
$ws.Range("A63").Select() | Out-Null
